# Add data for 2025-10-07
#
# Updates the "2025" year-to-date column (column L) with the latest figures
# across the Citywide Totals sheet, the By Neighborhood summary sheet, and the
# individual neighborhood detail sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 5200  # Aggravated Assault: 5183 -> 5200
$ws.Range("L3").Value = 5617  # Aggravated Battery: 5585 -> 5617
$ws.Range("L4").Value = 1363  # Criminal Sexual Assault: 1359 -> 1363
$ws.Range("L5").Value = 336  # Homicide: 332 -> 336
$ws.Range("L6").Value = 4672  # Robbery: 4653 -> 4672
$ws.Range("L7").Value = 17188  # Total: 17112 -> 17188

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 334  # Aggravated Assault: 329 -> 334
$ws.Range("L3").Value = 389  # Aggravated Battery: 386 -> 389
$ws.Range("L4").Value = 84  # Criminal Sexual Assault: 82 -> 84
$ws.Range("L5").Value = 40  # Homicide: 39 -> 40
$ws.Range("L6").Value = 292  # Robbery: 291 -> 292
$ws.Range("L7").Value = 1139  # Total: 1127 -> 1139

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 84  # Robbery: 83 -> 84
$ws.Range("L7").Value = 376  # Total: 375 -> 376

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 220  # Aggravated Assault: 219 -> 220
$ws.Range("L3").Value = 278  # Aggravated Battery: 276 -> 278
$ws.Range("L5").Value = 20  # Homicide: 19 -> 20
$ws.Range("L7").Value = 797  # Total: 793 -> 797

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L6").Value = 51  # Robbery: 50 -> 51
$ws.Range("L7").Value = 241  # Total: 240 -> 241

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 223  # Aggravated Battery: 222 -> 223
$ws.Range("L7").Value = 651  # Total: 650 -> 651

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L3").Value = 126  # Aggravated Battery: 125 -> 126
$ws.Range("L7").Value = 304  # Total: 303 -> 304

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 146  # Albany Park: 145 -> 146
$ws.Range("L5").Value = 62  # Armour Square: 61 -> 62
$ws.Range("L7").Value = 555  # Auburn Gresham: 553 -> 555
$ws.Range("L8").Value = 1139  # Austin: 1127 -> 1139
$ws.Range("L10").Value = 112  # Avondale: 111 -> 112
$ws.Range("L11").Value = 279  # Belmont Cragin: 276 -> 279
$ws.Range("L15").Value = 131  # Brighton Park: 130 -> 131
$ws.Range("L16").Value = 34  # Bucktown: 33 -> 34
$ws.Range("L19").Value = 462  # Chatham: 461 -> 462
$ws.Range("L20").Value = 425  # Chicago Lawn: 423 -> 425
$ws.Range("L21").Value = 54  # Chinatown: 53 -> 54
$ws.Range("L29").Value = 952  # Englewood: 945 -> 952
$ws.Range("L33").Value = 797  # Garfield Park: 793 -> 797
$ws.Range("L34").Value = 103  # Garfield Ridge: 102 -> 103
$ws.Range("L36").Value = 218  # Grand Boulevard: 217 -> 218
$ws.Range("L37").Value = 651  # Grand Crossing: 650 -> 651
$ws.Range("L42").Value = 562  # Humboldt Park: 561 -> 562
$ws.Range("L46").Value = 37  # Jefferson Park: 36 -> 37
$ws.Range("L48").Value = 221  # Lake View: 220 -> 221
$ws.Range("L52").Value = 349  # Little Village: 345 -> 349
$ws.Range("L54").Value = 369  # Loop: 367 -> 369
$ws.Range("L60").Value = 109  # Morgan Park: 108 -> 109
$ws.Range("L63").Value = 48  # NO NEIGHBORHOOD DATA: 50 -> 48
$ws.Range("L66").Value = 47  # North Center: 46 -> 47
$ws.Range("L67").Value = 589  # North Lawndale: 588 -> 589
$ws.Range("L70").Value = 29  # O'Hare: 28 -> 29
$ws.Range("L73").Value = 139  # Portage Park: 138 -> 139
$ws.Range("L75").Value = 62  # Pullman: 60 -> 62
$ws.Range("L76").Value = 266  # River North: 264 -> 266
$ws.Range("L77").Value = 115  # Riverdale: 113 -> 115
$ws.Range("L79").Value = 467  # Roseland: 462 -> 467
$ws.Range("L80").Value = 54  # Rush & Division: 53 -> 54
$ws.Range("L82").Value = 25  # Sheffield & DePaul: 24 -> 25
$ws.Range("L83").Value = 376  # South Chicago: 375 -> 376
$ws.Range("L84").Value = 169  # South Deering: 168 -> 169
$ws.Range("L85").Value = 866  # South Shore: 863 -> 866
$ws.Range("L86").Value = 119  # Streeterville: 120 -> 119
$ws.Range("L88").Value = 188  # United Center: 187 -> 188
$ws.Range("L91").Value = 230  # Washington Park: 229 -> 230
$ws.Range("L94").Value = 213  # West Loop: 212 -> 213
$ws.Range("L95").Value = 241  # West Pullman: 240 -> 241
$ws.Range("L96").Value = 194  # West Ridge: 193 -> 194
$ws.Range("L97").Value = 142  # West Town: 141 -> 142
$ws.Range("L99").Value = 304  # Woodlawn: 303 -> 304
$ws.Range("L101").Value = 17188  # Total: 17112 -> 17188

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L6").Value = 136  # Robbery: 135 -> 136
$ws.Range("L7").Value = 589  # Total: 588 -> 589

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 60  # Aggravated Battery: 59 -> 60
$ws.Range("L7").Value = 169  # Total: 168 -> 169

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 90  # Aggravated Battery: 88 -> 90
$ws.Range("L7").Value = 369  # Total: 367 -> 369

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 291  # Aggravated Assault: 290 -> 291
$ws.Range("L3").Value = 362  # Aggravated Battery: 357 -> 362
$ws.Range("L4").Value = 48  # Criminal Sexual Assault: 47 -> 48
$ws.Range("L7").Value = 952  # Total: 945 -> 952

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L4").Value = 44  # Criminal Sexual Assault: 43 -> 44
$ws.Range("L7").Value = 221  # Total: 220 -> 221

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L6").Value = 129  # Robbery: 128 -> 129
$ws.Range("L7").Value = 462  # Total: 461 -> 462

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("L2").Value = 46  # Aggravated Assault: 45 -> 46
$ws.Range("L3").Value = 33  # Aggravated Battery: 34 -> 33

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 49  # Aggravated Battery: 48 -> 49
$ws.Range("L6").Value = 125  # Robbery: 124 -> 125
$ws.Range("L7").Value = 266  # Total: 264 -> 266

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L4").Value = 46  # Criminal Sexual Assault: 45 -> 46
$ws.Range("L7").Value = 562  # Total: 561 -> 562

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("L3").Value = 30  # Aggravated Battery: 29 -> 30
$ws.Range("L7").Value = 112  # Total: 111 -> 112

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L2").Value = 11  # Aggravated Assault: 10 -> 11
$ws.Range("L7").Value = 37  # Total: 36 -> 37

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 59  # Aggravated Battery: 58 -> 59
$ws.Range("L7").Value = 194  # Total: 193 -> 194

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 104  # Aggravated Battery: 103 -> 104
$ws.Range("L7").Value = 230  # Total: 229 -> 230

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L6").Value = 28  # Robbery: 27 -> 28
$ws.Range("L7").Value = 54  # Total: 53 -> 54

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L3").Value = 158  # Aggravated Battery: 157 -> 158
$ws.Range("L6").Value = 116  # Robbery: 112 -> 116
$ws.Range("L7").Value = 467  # Total: 462 -> 467

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 133  # Aggravated Assault: 131 -> 133
$ws.Range("L7").Value = 425  # Total: 423 -> 425

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L6").Value = 55  # Robbery: 54 -> 55
$ws.Range("L7").Value = 218  # Total: 217 -> 218

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 187  # Aggravated Battery: 186 -> 187
$ws.Range("L5").Value = 14  # Homicide: 13 -> 14
$ws.Range("L7").Value = 555  # Total: 553 -> 555

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("L3").Value = 29  # Aggravated Battery: 28 -> 29
$ws.Range("L7").Value = 103  # Total: 102 -> 103

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L2").Value = 50  # Aggravated Assault: 49 -> 50
$ws.Range("L7").Value = 213  # Total: 212 -> 213

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L5").Value = 3  # Homicide: 2 -> 3
$ws.Range("L7").Value = 131  # Total: 130 -> 131

$ws = $wb.Worksheets.Item("North Center")
$ws.Range("L3").Value = 15  # Aggravated Battery: 14 -> 15
$ws.Range("L7").Value = 47  # Total: 46 -> 47

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 106  # Aggravated Assault: 105 -> 106
$ws.Range("L3").Value = 85  # Aggravated Battery: 83 -> 85
$ws.Range("L7").Value = 279  # Total: 276 -> 279

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L3").Value = 41  # Aggravated Battery: 40 -> 41
$ws.Range("L7").Value = 139  # Total: 138 -> 139

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L3").Value = 50  # Aggravated Battery: 49 -> 50
$ws.Range("L7").Value = 146  # Total: 145 -> 146

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("L3").Value = 33  # Aggravated Battery: 32 -> 33
$ws.Range("L7").Value = 142  # Total: 141 -> 142

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("L3").Value = 14  # Aggravated Battery: 13 -> 14
$ws.Range("L7").Value = 29  # Total: 28 -> 29

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L4").Value = 9  # Criminal Sexual Assault: 8 -> 9
$ws.Range("L7").Value = 188  # Total: 187 -> 188

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("L6").Value = 26  # Robbery: 25 -> 26
$ws.Range("L7").Value = 62  # Total: 61 -> 62

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 64  # Criminal Sexual Assault: 65 -> 64
$ws.Range("L7").Value = 119  # Total: 120 -> 119

$ws = $wb.Worksheets.Item("Pullman")
$ws.Range("L6").Value = 6  # Robbery: 4 -> 6
$ws.Range("L7").Value = 62  # Total: 60 -> 62

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 39  # Aggravated Battery: 38 -> 39
$ws.Range("L7").Value = 109  # Total: 108 -> 109

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 354  # Aggravated Battery: 352 -> 354
$ws.Range("L6").Value = 179  # Robbery: 178 -> 179
$ws.Range("L7").Value = 866  # Total: 863 -> 866

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("L6").Value = 8  # Robbery: 7 -> 8
$ws.Range("L7").Value = 25  # Total: 24 -> 25

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L2").Value = 40  # Aggravated Assault: 39 -> 40
$ws.Range("L3").Value = 38  # Aggravated Battery: 37 -> 38
$ws.Range("L7").Value = 115  # Total: 113 -> 115

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 28  # Robbery: 27 -> 28
$ws.Range("L7").Value = 54  # Total: 53 -> 54

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L2").Value = 114  # Aggravated Assault: 111 -> 114
$ws.Range("L3").Value = 112  # Aggravated Battery: 111 -> 112
$ws.Range("L7").Value = 349  # Total: 345 -> 349

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("L6").Value = 22  # Robbery: 21 -> 22
$ws.Range("L7").Value = 34  # Total: 33 -> 34
